$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indexes: A=1 time, B=2 type, C=3 user_name, D=4 pc_name, E=5 win_title,
# F=6 win_urlpath, G=7 win_hwnd, H=8 win_class, I=9 app_path, J=10 capimg,
# K=11 explanation, L=12 error_type, M=13 error_content

# Row 2
$ws.Cells.Item(2, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(2, 10).Value = "bdot20240415_141954/1.png"
$ws.Cells.Item(2, 11).Value = "「スタート」ボタンをクリックする"

# Row 3
$ws.Cells.Item(3, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(3, 10).Value = "bdot20240415_141954/2.png"
$ws.Cells.Item(3, 11).Value = "メニューから「設定」アイコンをクリックする"

# Row 4
$ws.Cells.Item(4, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(4, 10).Value = "bdot20240415_141954/3.png"
$ws.Cells.Item(4, 11).Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# Row 5
$ws.Cells.Item(5, 2).Value = "error"
$ws.Cells.Item(5, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(5, 10).Value = "bdot20240415_141954/4.png"
$ws.Cells.Item(5, 11).Value = "0x80240fff エラー"
$ws.Cells.Item(5, 12).Value = "Error W"
$ws.Cells.Item(5, 13).Value = " エラーの Windows"

# Row 6
$ws.Cells.Item(6, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(6, 10).Value = "bdot20240415_141954/5.png"
$ws.Cells.Item(6, 11).Value = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"

# Row 7
$ws.Cells.Item(7, 2).Value = "operation"
$ws.Cells.Item(7, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(7, 10).Value = "bdot20240415_141954/5.png"
$ws.Cells.Item(7, 11).Value = "メニューからターミナル(管理者)をクリックする"
$ws.Cells.Item(7, 12).Value = ""
$ws.Cells.Item(7, 13).Value = ""

# Row 8
$ws.Cells.Item(8, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(8, 10).Value = "bdot20240415_141954/6.png"
$ws.Cells.Item(8, 11).Value = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"

# Row 9
$ws.Cells.Item(9, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(9, 10).Value = "bdot20240415_141954/7.png"
$ws.Cells.Item(9, 11).Value = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"

# Row 10
$ws.Cells.Item(10, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(10, 10).Value = "bdot20240415_141954/8.png"
$ws.Cells.Item(10, 11).Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"

# Row 11
$ws.Cells.Item(11, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(11, 10).Value = "bdot20240415_141954/9.png"
$ws.Cells.Item(11, 11).Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"

# Row 12
$ws.Cells.Item(12, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(12, 10).Value = "bdot20240415_141954/10.png"
$ws.Cells.Item(12, 11).Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"

# Row 13
$ws.Cells.Item(13, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(13, 10).Value = "bdot20240415_141954/1.png"
$ws.Cells.Item(13, 11).Value = "「スタート」ボタンをクリックする"

# Row 14
$ws.Cells.Item(14, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(14, 10).Value = "bdot20240415_141954/2.png"
$ws.Cells.Item(14, 11).Value = "メニューから「設定」アイコンをクリックする"

# Row 15
$ws.Cells.Item(15, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(15, 10).Value = "bdot20240415_141954/3.png"
$ws.Cells.Item(15, 11).Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# Row 16
$ws.Cells.Item(16, 3).Value = "Nanami Ishikawa"
$ws.Cells.Item(16, 10).Value = "bdot20240415_141954/11.png"
$ws.Cells.Item(16, 11).Value = "「更新プログラムのチェック」ボタンをクリックする"
